# Auto-generated Word COM-interop script implementing the assignment-nine
# Requirements.docx edit described by the commit "update to assn 9".
#
# Strategy: the edit restructures the U01/U02 functional-requirement
# paragraphs (splitting the "U01"/"U02" labels into separate runs and
# moving the "U01-2"/"U02-2" detail paragraphs up next to them), rewords
# the "General non-functional requirements" sentence across several runs,
# trims the "U01-1"/"U02-1" sentences (dropping their labels & one
# lastRenderedPageBreak, adding "within 1.5 seconds" timing), and moves a
# lastRenderedPageBreak marker from the "Interface for Repo-Group trend:"
# heading down to the "Interface with User Trend:" picture run.
#
# Because this reshuffles/splits so many runs and paragraphs at once, the
# most reliable way to express it through the Word OM is with
# Range.InsertXML against the exact paragraph span being changed (using a
# minimal WordProcessingML "pkg:package" payload), rather than chaining
# many small Find/Replace or Range.Text edits. Paragraph indices (not
# Find-result Range math, which proved unreliable for Next()/MoveEnd in
# this host) are used to pin down the exact Range boundaries.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, $substring) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$substring*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------
# Step 1: locate the paragraph range "U01: Application ... select group."
# through the trailing empty paragraph right after "U02-2: ... user name."
# (this is the whole functional/non-functional requirements block that
# gets reshuffled) and replace it in one shot.
# ---------------------------------------------------------------
$idxStart = Find-ParagraphIndex $d "U01: Application will display repo group commit activity for select group."
$idxU022 = Find-ParagraphIndex $d "U02-2: Line chart will have number of commits as the vertical axis and date as the horizontal axis and the line will be label with the user name."

$pStart = $d.Paragraphs.Item($idxStart)
$pU022 = $d.Paragraphs.Item($idxU022)
$pTrailingEmpty = $pU022.Next()

$rng1 = $d.Range($pStart.Range.Start, $pTrailingEmpty.Range.End)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D47BCE" w:rsidRDefault="00D47BCE"><w:r><w:tab/><w:t>U01</w:t></w:r><w:r><w:t>-1</w:t></w:r><w:r><w:t>: Application will display repo group commit activity for select group.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>U01-2: Line chart will have number of commits as the vertical axis and date as the horizontal axis and the line will be label with the repo group name.</w:t></w:r></w:p><w:p w:rsidR="00D47BCE" w:rsidRDefault="00D47BCE"><w:r><w:tab/><w:t>U02</w:t></w:r><w:r><w:t>-1</w:t></w:r><w:r><w:t>: Application will display results of commit activity of selected user in selected repo group.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>U02-2: Line chart will have number of commits as the vertical axis and date as the horizontal axis and the line will be label with the user name.</w:t></w:r></w:p><w:p/><w:p w:rsidR="00117C41" w:rsidRDefault="00050572" w:rsidP="00117C41"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">4.2 </w:t></w:r><w:r w:rsidR="00117C41" w:rsidRPr="00D47BCE"><w:rPr><w:b/></w:rPr><w:t>Non-function</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>al</w:t></w:r><w:r w:rsidR="00117C41" w:rsidRPr="00D47BCE"><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Requirements:</w:t></w:r></w:p><w:p w:rsidR="00050572" w:rsidRPr="00050572" w:rsidRDefault="00050572" w:rsidP="00356007"><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t>General non-functional requirem</w:t></w:r><w:r><w:t>ents: System will display correct</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">and validated </w:t></w:r><w:r><w:t>data for all repo group and user options</w:t></w:r><w:r><w:t xml:space="preserve"> within 1.5 seconds</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00117C41" w:rsidRDefault="00D47BCE" w:rsidP="00050572"><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t>Application will have line chart displaying results</w:t></w:r><w:r><w:t xml:space="preserve"> within 1.5 seconds</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p w:rsidR="00050572" w:rsidRDefault="00050572" w:rsidP="00050572"><w:pPr><w:ind w:firstLine="720"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng1.InsertXML($xml1)

# ---------------------------------------------------------------
# Step 2: drop the lastRenderedPageBreak before "Interface for Repo-Group trend:"
# ---------------------------------------------------------------
$idx2 = Find-ParagraphIndex $d "Interface for Repo-Group trend:"
$p2 = $d.Paragraphs.Item($idx2)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00160204" w:rsidRPr="00356007" w:rsidRDefault="00160204" w:rsidP="00050572"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r w:rsidRPr="00356007"><w:rPr><w:b/></w:rPr><w:t>Interface for Repo-Group trend:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------
# Step 3: add a lastRenderedPageBreak before the drawing that follows
# "Interface with User Trend:" (the InterfaceUserExample.png picture run).
# This paragraph happens to be the very last paragraph in the document
# body (immediately followed by the sectPr), and InsertXML on the last
# paragraph's Range leaves a stray trailing empty paragraph behind (the
# payload's own implicit final paragraph mark doesn't merge with the
# document's existing one) -- so an extra cleanup step removes it.
# ---------------------------------------------------------------
$idx3 = Find-ParagraphIndex $d "Interface with User Trend:"
$headingPara = $d.Paragraphs.Item($idx3)
$picPara = $headingPara.Next()
$countBefore = $d.Paragraphs.Count
$xml3 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w:rsidR="00160204" w:rsidRPr="00160204" w:rsidRDefault="00160204" w:rsidP="00160204"><w:pPr><w:tabs><w:tab w:val="left" w:pos="2160"/></w:tabs></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="4320938" cy="3501529"/><wp:effectExtent l="19050" t="0" r="3412" b="0"/><wp:docPr id="2" name="Picture 1" descr="InterfaceUserExample.png"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="InterfaceUserExample.png"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId6" cstate="print"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="4319904" cy="3500691"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$picPara.Range.InsertXML($xml3)

if ($d.Paragraphs.Count -gt $countBefore) {
    $newLast = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newSecondLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    if ($newLast.Range.Start -eq ($newLast.Range.End - 1)) {
        $cleanupRng = $d.Range($newSecondLast.Range.End - 1, $newLast.Range.End - 1)
        $cleanupRng.Delete()
    }
}

Write-Output "edit complete"
